$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COMO_History")

# Update the DSC1/DSC2 note text (rows 2 and 3, column I) with revised descriptions
$ws.Range("I2").Value = "DSC1 is the original discharge XS, located 1 meter downstream of S2 and the temporary staff gauge (Gauge1) in the tail of a pool.  DSC1 was never surveyed but Gauge1 was during geomorphology.  "
$ws.Range("I3").Value = "DSC2 is the current DSC XS, it is located 30 meters upstream of S2 and Gauge2 in a riffle so a different hydrologic unit as the Level TROLL and staff gauge.  DSC2 and Gauge2 were surveyed during geomorphology.  Hannah says it's possible that they began collecting DSC measurements at DSC2 before Gauge2 was installed in 2016.  Hannah thinks that the first measurement associated with DSC2 may have been 6/13/16.  "

# DSC2's "Installed" date is now known - fill in the date instead of the old "?" placeholder
$ws.Range("D3").Value = 42675
$ws.Range("D3").NumberFormat = "mmm-yy"

# Make COMO_History the active/selected sheet with F10 selected
$ws.Activate()
$ws.Range("F10").Select()
